$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old table was a 4x4 block (Source | A | FFR | LF  with A/FFR/LF Lag rows).
# The FFR column and FFR Lag row are removed, shrinking the table to 3x3, and
# the numeric entries are refreshed with new Crisis/Credit-Allocation figures.

# Drop the surplus FFR column (D) and the surplus FFR Lag row (row 4) entirely
# so the used range shrinks back down to A1:C3.
$ws.Range("D1:D4").Clear()
$ws.Range("A4:C4").Clear()

# Header row: "FFR" -> "LF" (A1/B1 stay the same: Source / A).
$ws.Range("C1").Value = "LF"

# Row labels: second row stays "A Lag"; the old "FFR Lag" row label is now
# repurposed as "LF Lag" on row 3.
$ws.Range("A3").Value = "LF Lag"

# Refreshed coefficients.
# "0.27" reads as a plain number to the smart-input parser, which would store
# it as a numeric cell and stamp a new text-quoted style on it; the source
# workbook keeps every coefficient (even purely numeric-looking ones) as a
# plain shared-string with no special formatting, so round-trip it through a
# text formula + paste-values to land it as an ordinary string cell.
$ws.Range("B2").Formula = '="0.27"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("C2").Value = "-8.77*"
$ws.Range("B3").Value = "-0.11*"
$ws.Range("C3").Value = "2.45***"
